$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Certificados, cursos, badges")

# Duplicate the formatting of the last existing row (73) onto the two new rows.
$ws.Range("B73:I73").Copy()
$ws.Range("B74:I75").PasteSpecial(-4122)

# Row 74 - Big Data Real-Time Analytics com Python e Spark
$ws.Cells.Item(74, 2).Value = "Data Science Academy"
$ws.Cells.Item(74, 3).Value = "Big Data Real-Time Analytics com Python e Spark"
$ws.Cells.Item(74, 4).Value = 72
$ws.Cells.Item(74, 5).Value = 45536
$ws.Hyperlinks.Add($ws.Cells.Item(74, 6), "https://mycourse.app/CcubvCKzPXCtRwHf7")
$ws.Cells.Item(74, 7).Value = "Ok"
$ws.Cells.Item(74, 8).Value = "Ok"
$ws.Cells.Item(74, 9).Value = 45536

# Row 75 - Visualização de Dados e Design de Dashboards
$ws.Cells.Item(75, 2).Value = "Data Science Academy"
$ws.Cells.Item(75, 3).Value = "Visualização de Dados e Design de Dashboards"
$ws.Cells.Item(75, 4).Value = 54
$ws.Cells.Item(75, 5).Value = 45536
$ws.Hyperlinks.Add($ws.Cells.Item(75, 6), "https://mycourse.app/vyPDKbNYCDmw8EPP8")
$ws.Cells.Item(75, 7).Value = "Ok"
$ws.Cells.Item(75, 8).Value = "Ok"
$ws.Cells.Item(75, 9).Value = 45536

# Re-apply row 73's formatting on top so the hyperlink cells (F74/F75) keep the
# workbook's custom "Hiperlink" look (9pt underlined + border) instead of the
# generic one Hyperlinks.Add uses.
$ws.Range("B73:I73").Copy()
$ws.Range("B74:I75").PasteSpecial(-4122)

$ws.Range("A2").Select()
